$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tables")
$ws.Rows("99:99").Insert()
